$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. "geo" row clarity: the [float,float,float?] description becomes
#    [float,float] (M6, next to the "geo" row in the side table).
# ---------------------------------------------------------------------------
$ws.Range("M6").Value = "[float,float]"

# ---------------------------------------------------------------------------
# 2. Insert a new "reserved" bit-pattern row before the SEMANTIC row (26),
#    shifting SEMANTIC (and the two trailing blank spacer rows) down by one.
#    Excel correctly re-adjusts the self-referential formulas and inherits
#    the formatting of the row above, just like a manual row insert.
# ---------------------------------------------------------------------------
$ws.Rows.Item(26).Insert()

# ---------------------------------------------------------------------------
# 3. Re-label two of the "reserved" rows as the newly introduced
#    "list-of-lists" / "list-of-maps" types, and touch up the bit patterns
#    for rows 22-26 so the reserved/new rows keep a contiguous sequence.
# ---------------------------------------------------------------------------

# Row 22: reserved, bits -> 1,1,0 (48)
$ws.Range("F22").Value = 1
$ws.Range("G22").Value = 1
$ws.Range("H22").Value = 0

# Row 23: reserved, bits -> 0,0,1 (64)
$ws.Range("F23").Value = 0
$ws.Range("G23").Value = 0
$ws.Range("H23").Value = 1

# Row 24: becomes "list-of-lists", bits -> 1,0,1 (80)
$ws.Range("A24").Value = "list-of-lists"
$ws.Range("F24").Value = 1
$ws.Range("G24").Value = 0
$ws.Range("H24").Value = 1

# Row 25: becomes "list-of-maps", bits -> 0,1,1 (96)
$ws.Range("A25").Value = "list-of-maps"
$ws.Range("F25").Value = 0
$ws.Range("G25").Value = 1
$ws.Range("H25").Value = 1

# Row 26 (newly inserted row): reserved, bits -> 1,1,1 (112)
$ws.Range("A26").Value = "reserved"
$ws.Range("A26").Style = $ws.Range("A25").Style
$ws.Range("F26").Value = 1
$ws.Range("G26").Value = 1
$ws.Range("H26").Value = 1
$ws.Range("J26").Formula = "=(B26*B3)+(C26*C3)+(D26*D3)+(E26*E3)+(F26*F3)+(G26*G3)+(H26*H3)+(I26*I3)"
